$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C:F, rows 2-21 (biomass, product, substrate, volume)
$data = @(
    @(0.1665997958244148, 56.43332382052655, 0.153759167860843, 0.5928920083319424),
    @(0.3657147854415364, 41.45415687577069, 2.456545809447455, 0.7680803126325425),
    @(0.5398665447771468, 33.62353661399179, 3.518019616177737, 0.9739507866271381),
    @(0.7223430829397475, 26.21012337465516, 3.966217411450854, 1.306251845543919),
    @(1.26110052964417, 23.65686600649336, 3.561087035725909, 1.515346477858623),
    @(2.072541300615169, 21.36152468615415, 2.38198148263037, 1.670874335642752),
    @(3.070328327164224, 19.60653101987885, 0.8570223710457577, 1.88366409190616),
    @(3.601914334421808, 16.85870156828176, 0.237852768163538, 2.143990387381906),
    @(3.701880253661649, 14.93160137600482, 0.2306129846327113, 2.250844921357771),
    @(3.772264407191295, 14.2067572290145, 0.1277857790922204, 2.509038325505864),
    @(3.86461863085781, 13.49649419094325, 0.180115877911402, 2.783234101975589),
    @(3.969546192087241, 13.13395011661649, 0.06363411700330768, 2.838674410909438),
    @(4.008608278221033, 12.01657681006507, 0.1109631015419993, 3.166501999639347),
    @(4.166151155780338, 10.27239061750221, 0.09094663140392872, 3.300747313374659),
    @(4.08189404620815, 11.09237726061986, 0.1223274968641734, 3.666612573073557),
    @(4.323848634415815, 9.613567789615315, 0.0739916204664739, 3.767131880376319),
    @(4.291984100163914, 10.06414044836635, 0.07377692416886689, 4.170961444635592),
    @(4.324296546215024, 9.257749308872018, 0.09913714621632187, 4.057220001967738),
    @(4.523842889773126, 8.055541131311198, 0.1056686815460239, 4.321908806023299),
    @(4.484353700269675, 7.740824033395551, 0.05744272051531361, 4.649154369382726)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}
